$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.907.15"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "2.485.30"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.23%  "
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.23"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.126"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("E14").Value = "  -2.37%  "
$ws.Range("D15").Value = "2.872.78"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "2.487.63"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.827"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "47.790.87"
$ws.Range("E18").Value = "  -0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.54%  "
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "281.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("E25").Value = "  -3.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("E29").Value = "  -5.57%  "
$ws.Range("E30").Value = "  -3.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.36%  "
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0766"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.00%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -3.93%  "
$ws.Range("E39").Value = "  -3.63%  "
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("E44").Value = "  -0.63%  "
$ws.Range("D45").Value = "1.983.85"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.20%  "
